# Apply updated crypto price/volume figures (commit: "Updated cryptos list on Sat Jul  6 07:42:28 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values that look like plain numbers get a leading apostrophe so Excel
# stores them as literal text (preserving trailing zeros / exact formatting), matching the
# original text-typed cells instead of being auto-converted to numeric values.

$ws.Range("D2").Value = '56.615.63'
$ws.Range("E2").Value = '  +4.68%  '
$ws.Range("D3").Value = '3.014.27'
$ws.Range("E3").Value = '  +5.54%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''509.23'
$ws.Range("E5").Value = '  +8.74%  '
$ws.Range("D6").Value = '''137.23'
$ws.Range("E6").Value = '  +9.36%  '
$ws.Range("E8").Value = '  +8.41%  '
$ws.Range("D9").Value = '''7.62'
$ws.Range("E9").Value = '  +15.50%  '
$ws.Range("E10").Value = '  +13.81%  '
$ws.Range("D11").Value = '''0.355'
$ws.Range("E11").Value = '  +8.06%  '
$ws.Range("E12").Value = '  +4.80%  '
$ws.Range("D13").Value = '3.526.35'
$ws.Range("E13").Value = '  +5.39%  '
$ws.Range("D14").Value = '''25.69'
$ws.Range("E14").Value = '  +11.09%  '
$ws.Range("D15").Value = '''0.0000155'
$ws.Range("E15").Value = '  +16.25%  '
$ws.Range("D16").Value = '56.619.08'
$ws.Range("E16").Value = '  +4.67%  '
$ws.Range("D17").Value = '3.008.58'
$ws.Range("E17").Value = '  +5.08%  '
$ws.Range("E18").Value = '  +9.85%  '
$ws.Range("D19").Value = '''12.54'
$ws.Range("E19").Value = '  +10.16%  '
$ws.Range("E20").Value = '  +12.53%  '
$ws.Range("D21").Value = '''327.86'
$ws.Range("E21").Value = '  +12.22%  '
$ws.Range("D23").Value = '''0.479'
$ws.Range("E23").Value = '  +8.66%  '
$ws.Range("D24").Value = '''62.54'
$ws.Range("E24").Value = '  +7.09%  '
$ws.Range("E25").Value = '  +10.40%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '0.0₃0925'
$ws.Range("E27").Value = '  +15.09%  '
$ws.Range("D28").Value = '''6.58'
$ws.Range("E28").Value = '  +7.30%  '
$ws.Range("D29").Value = '''7.00'
$ws.Range("E29").Value = '  +13.68%  '
$ws.Range("E30").Value = '  +11.49%  '
$ws.Range("E31").Value = '  +10.05%  '
$ws.Range("D32").Value = '''20.70'
$ws.Range("E32").Value = '  +10.41%  '
$ws.Range("D33").Value = '''156.14'
$ws.Range("E33").Value = '  +16.40%  '
$ws.Range("D34").Value = '''4.53'
$ws.Range("E34").Value = '  +8.23%  '
$ws.Range("D35").Value = '''5.63'
$ws.Range("E35").Value = '  +4.37%  '
$ws.Range("D36").Value = '''1.28'
$ws.Range("E36").Value = '  +5.16%  '
$ws.Range("D37").Value = '''0.0676'
$ws.Range("E37").Value = '  +9.92%  '
$ws.Range("D38").Value = '''23.83'
$ws.Range("E38").Value = '  +3.83%  '
$ws.Range("D39").Value = '3.046.00'
$ws.Range("E39").Value = '  +5.79%  '
$ws.Range("D40").Value = '''36.61'
$ws.Range("E40").Value = '  +5.43%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  +8.16%  '
$ws.Range("D43").Value = '2.268.47'
$ws.Range("E43").Value = '  +11.54%  '
$ws.Range("D44").Value = '''1.01'
$ws.Range("E44").Value = '  +6.57%  '
$ws.Range("E45").Value = '  +7.55%  '
$ws.Range("D46").Value = '''3.62'
$ws.Range("E46").Value = '  +7.57%  '
$ws.Range("D47").Value = '''2.01'
$ws.Range("E47").Value = '  +25.92%  '
$ws.Range("D48").Value = '''0.0237'
$ws.Range("E48").Value = '  +11.99%  '
$ws.Range("E49").Value = '  +8.74%  '
$ws.Range("D50").Value = '''19.23'
$ws.Range("E50").Value = '  +8.07%  '
$ws.Range("E51").Value = '  +11.51%  '
